$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and 1h volume change (E) columns with latest scraped values.
# NumberFormat is forced to text ("@") before assigning so values such as "40.187.26"
# or "0.0719" are not auto-converted to numbers/dates by Excel, then the style is reset
# back to "Normal" so no stray formatting is left behind on the cell.

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "40.187.26"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  +2.72%  "
$cell.Style = "Normal"

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.247.04"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  +0.59%  "
$cell.Style = "Normal"

# Row 4
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  -0.10%  "
$cell.Style = "Normal"

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "295.54"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -0.72%  "
$cell.Style = "Normal"

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "87.45"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  +7.92%  "
$cell.Style = "Normal"

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.515"
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  +1.05%  "
$cell.Style = "Normal"

# Row 8
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  +0.01%  "
$cell.Style = "Normal"

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.475"
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  +2.89%  "
$cell.Style = "Normal"

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "31.20"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  +11.11%  "
$cell.Style = "Normal"

# Row 11
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  +3.59%  "
$cell.Style = "Normal"

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "47.17"
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  +2.13%  "
$cell.Style = "Normal"

# Row 13
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  +0.56%  "
$cell.Style = "Normal"

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.47"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  +5.56%  "
$cell.Style = "Normal"

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "2.590.35"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  +0.47%  "
$cell.Style = "Normal"

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "14.30"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  +1.48%  "
$cell.Style = "Normal"

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.247.88"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  +0.28%  "
$cell.Style = "Normal"

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.736"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  +2.56%  "
$cell.Style = "Normal"

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "40.096.27"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  +2.66%  "
$cell.Style = "Normal"

# Row 20
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  +3.71%  "
$cell.Style = "Normal"

# Row 21
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  +2.39%  "
$cell.Style = "Normal"

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "10.69"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  +7.56%  "
$cell.Style = "Normal"

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "65.86"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  +0.95%  "
$cell.Style = "Normal"

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "236.87"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  +4.73%  "
$cell.Style = "Normal"

# Row 25
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  +0.02%  "
$cell.Style = "Normal"

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.48"
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  +3.45%  "
$cell.Style = "Normal"

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "23.22"
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  +3.99%  "
$cell.Style = "Normal"

# Row 29
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  +5.00%  "
$cell.Style = "Normal"

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "9.30"
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  +4.11%  "
$cell.Style = "Normal"

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "34.23"
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  +8.80%  "
$cell.Style = "Normal"

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "153.34"
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  +2.66%  "
$cell.Style = "Normal"

# Row 33
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -0.14%  "
$cell.Style = "Normal"

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "4.92"
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  +2.67%  "
$cell.Style = "Normal"

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.0719"
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  +4.87%  "
$cell.Style = "Normal"

# Row 36
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  +2.40%  "
$cell.Style = "Normal"

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "16.77"
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  +13.84%  "
$cell.Style = "Normal"

# Row 38
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  +5.73%  "
$cell.Style = "Normal"

# Row 39
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  +2.38%  "
$cell.Style = "Normal"

# Row 40
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  +2.01%  "
$cell.Style = "Normal"

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.70"
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  +5.48%  "
$cell.Style = "Normal"

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "3.84"
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  +4.93%  "
$cell.Style = "Normal"

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.997.87"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  +4.74%  "
$cell.Style = "Normal"

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.23"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  +7.48%  "
$cell.Style = "Normal"

# Row 45
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  +6.87%  "
$cell.Style = "Normal"

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "10.03"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  +11.61%  "
$cell.Style = "Normal"

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "16.50"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  +0.61%  "
$cell.Style = "Normal"

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.60"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  +2.22%  "
$cell.Style = "Normal"

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "2.460.48"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  +0.77%  "
$cell.Style = "Normal"

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "71.55"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  +6.64%  "
$cell.Style = "Normal"

# Row 51
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  +14.53%  "
$cell.Style = "Normal"
